$wb = $excel.ActiveWorkbook

# --- Sheet "List": rename the accountId header/value to account ---
$list = $wb.Worksheets.Item("List")
$list.Range("E1").Value = "`${msg.getProperty('savedSearch_account')}"
$list.Range("E2").Value = "`${printer.print(savedSearch.account)}"

# --- Sheet "Search": add a new row 7 mirroring the account column ---
$search = $wb.Worksheets.Item("Search")
$search.Range("A7").Value = "`${msg.getProperty('savedSearch_account')}"
$search.Range("B7").Value = "`${account}"
